# Insert a new weekly data row before row 32 (shifts existing rows 32-180 down to 33-181)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with the new weekly record
$ws.Range("A32").Value = 4
$ws.Range("B32").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C32").Value = "Los Lagos"
$ws.Range("D32").Value = 44560
$ws.Range("E32").Value = 10
$ws.Range("F32").Value = 100112017
$ws.Range("G32").Value = "Apio"
$ws.Range("H32").Value = "Americana (o)"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 20
$ws.Range("K32").Value = 12000
$ws.Range("L32").Value = 12500
$ws.Range("M32").Value = 12250
$ws.Range("N32").Value = "$/docena de matas"
$ws.Range("O32").Value = "Región de Coquimbo"
$ws.Range("P32").Value = 2042
$ws.Range("Q32").Value = 6
$ws.Range("R32").Value = "Hortaliza"
